$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.611.13"
$ws.Range("E2").Value = "  +5.32%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.243.42"
$ws.Range("E3").Value = "  +2.66%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'577.79"
$ws.Range("E5").Value = "  +2.85%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'179.85"
$ws.Range("E6").Value = "  +6.77%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -2.19%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.240.25"
$ws.Range("E9").Value = "  +2.57%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +4.57%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'6.79"
$ws.Range("E11").Value = "  +3.79%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +4.89%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.808.08"
$ws.Range("E13").Value = "  +2.88%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.15%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "'27.99"
$ws.Range("E15").Value = "  +3.24%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "67.481.68"
$ws.Range("E16").Value = "  +5.16%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "'0.0000167"
$ws.Range("E17").Value = "  +2.91%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.240.82"
$ws.Range("E18").Value = "  +2.74%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.48%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +3.78%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'374.45"
$ws.Range("E21").Value = "  +6.45%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +5.64%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.54%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'71.23"
$ws.Range("E24").Value = "  +4.64%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  +1.95%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +3.86%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  +1.36%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +3.15%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.05%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +4.86%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "'5.61"
$ws.Range("E31").Value = "  +2.99%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "'22.61"
$ws.Range("E32").Value = "  +3.26%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  -0.10%  "

# Row 34 - Fetch.AI
$ws.Range("D34").Value = "'1.27"
$ws.Range("E34").Value = "  +6.65%  "

# Row 35 - Aptos
$ws.Range("D35").Value = "'6.84"
$ws.Range("E35").Value = "  +3.87%  "

# Row 36 - Monero
$ws.Range("D36").Value = "'164.62"
$ws.Range("E36").Value = "  +6.63%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "'1.49"
$ws.Range("E37").Value = "  +4.88%  "

# Row 38 - Mantle
$ws.Range("D38").Value = "'0.863"
$ws.Range("E38").Value = "  +5.78%  "

# Row 39 - Stacks
$ws.Range("D39").Value = "'1.86"
$ws.Range("E39").Value = "  +10.27%  "

# Row 40 - RenderToken
$ws.Range("D40").Value = "'6.88"
$ws.Range("E40").Value = "  +15.02%  "

# Row 41 - EnergySwap
$ws.Range("D41").Value = "'26.72"
$ws.Range("E41").Value = "  +1.11%  "

# Row 42 - was Bittensor, now dogwifhat (swap with row 43)
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.59"
$ws.Range("E42").Value = "  +5.51%  "

# Row 43 - was dogwifhat, now Bittensor (swap with row 42)
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'362.76"
$ws.Range("E43").Value = "  +13.25%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  +6.03%  "

# Row 45 - Maker
$ws.Range("D45").Value = "2.700.36"
$ws.Range("E45").Value = "  +2.34%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "'25.82"
$ws.Range("E46").Value = "  +9.36%  "

# Row 47 - OKB
$ws.Range("D47").Value = "'40.45"
$ws.Range("E47").Value = "  +2.77%  "

# Row 48 - Hedera
$ws.Range("D48").Value = "'0.0673"
$ws.Range("E48").Value = "  +3.80%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +2.95%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  +1.14%  "

# Row 51 - ONDO
$ws.Range("E51").Value = "  +6.73%  "
